$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" cells: force text storage (values look numeric, e.g. "579.46")
# so they must be written as text, matching the source data which uses inline/shared
# strings rather than numeric cells (some prices use "." as a thousands separator,
# e.g. "66.339.18", so the whole column is textual).
$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "66.339.18"
$r.Style = "Normal"
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "3.085.28"
$r.Style = "Normal"
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "579.46"
$r.Style = "Normal"
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "167.85"
$r.Style = "Normal"
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.Style = "Normal"
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "3.088.80"
$r.Style = "Normal"
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "6.67"
$r.Style = "Normal"
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.153"
$r.Style = "Normal"
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "36.76"
$r.Style = "Normal"
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "3.583.81"
$r.Style = "Normal"
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "66.330.27"
$r.Style = "Normal"
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "3.072.06"
$r.Style = "Normal"
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "16.10"
$r.Style = "Normal"
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "464.78"
$r.Style = "Normal"
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "0.713"
$r.Style = "Normal"
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "7.45"
$r.Style = "Normal"
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "83.14"
$r.Style = "Normal"
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "2.28"
$r.Style = "Normal"
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "10.12"
$r.Style = "Normal"
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "28.36"
$r.Style = "Normal"
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.Style = "Normal"
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "5.88"
$r.Style = "Normal"
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "48.73"
$r.Style = "Normal"
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "50.10"
$r.Style = "Normal"
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.314"
$r.Style = "Normal"
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "2.91"
$r.Style = "Normal"
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "0.0361"
$r.Style = "Normal"
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "382.15"
$r.Style = "Normal"
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "2.768.73"
$r.Style = "Normal"
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "134.19"
$r.Style = "Normal"
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "24.51"
$r.Style = "Normal"

# E-column "Volume(1h)" cells: plain text percentages (already non-numeric due to
# padding/percent sign), safe to assign directly.
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("E3").Value = "  +4.20%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("E6").Value = "  +3.96%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  +4.41%  "
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("E11").Value = "  +1.04%  "
$ws.Range("E12").Value = "  +5.94%  "
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("E14").Value = "  +7.60%  "
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("E16").Value = "  +3.81%  "
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("E18").Value = "  +4.18%  "
$ws.Range("E19").Value = "  +3.44%  "
$ws.Range("E21").Value = "  +3.83%  "
$ws.Range("E22").Value = "  +5.68%  "
$ws.Range("E23").Value = "  +3.83%  "
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("E25").Value = "  +5.22%  "
$ws.Range("E26").Value = "  +2.70%  "
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("E31").Value = "  +3.24%  "
$ws.Range("E32").Value = "  +3.08%  "
$ws.Range("E33").Value = "  +4.40%  "
$ws.Range("E34").Value = "  +5.44%  "
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("E36").Value = "  +2.23%  "
$ws.Range("E37").Value = "  +3.07%  "
$ws.Range("E38").Value = "  +11.46%  "
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("E40").Value = "  +4.47%  "
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("E42").Value = "  +3.04%  "
$ws.Range("E43").Value = "  +2.90%  "
$ws.Range("E44").Value = "  +3.97%  "
$ws.Range("E45").Value = "  +2.42%  "
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("E48").Value = "  +2.75%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  +6.18%  "
$ws.Range("E51").Value = "  +4.77%  "
